# Auto-generated edit script applying numeric updates per the commit diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 33367.25
$ws.Range("I21").Value = 70019
$ws.Range("J21").Value = 21150
$ws.Range("K21").Value = 70019
$ws.Range("L21").Value = 21150
$ws.Range("M21").Value = -69551
$ws.Range("N21").Value = -22086

$ws.Range("H23").Value = 33367.25
$ws.Range("I23").Value = 70019
$ws.Range("J23").Value = 21150
$ws.Range("K23").Value = 70019
$ws.Range("L23").Value = 21150
$ws.Range("M23").Value = -69785
$ws.Range("N23").Value = -21618

$ws.Range("H32").Value = 3999.5
$ws.Range("I32").Value = 3999.5
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 3999.5
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -3673.5
$ws.Range("N32").ClearContents()

$ws.Range("H41").Value = 375.30768
$ws.Range("I41").Value = 209.875
$ws.Range("K41").Value = 209.875
$ws.Range("M41").Value = 230.125

$ws.Range("H43").Value = 994.5
$ws.Range("I43").Value = 780
$ws.Range("K43").Value = 780
$ws.Range("M43").Value = -711

$ws.Range("H51").Value = 2872
$ws.Range("J51").Value = 2965
$ws.Range("L51").Value = 2965
$ws.Range("N51").Value = -3933

$ws.Range("H76").Value = 2874.8293
$ws.Range("I76").Value = 2582.4333
$ws.Range("J76").Value = 3672.2727
$ws.Range("K76").Value = 2582.4333
$ws.Range("L76").Value = 3672.2727
$ws.Range("M76").Value = -2267.4333
$ws.Range("N76").Value = -4302.2727

$ws.Range("H79").Value = 2874.8293
$ws.Range("I79").Value = 2582.4333
$ws.Range("J79").Value = 3672.2727
$ws.Range("K79").Value = 2582.4333
$ws.Range("L79").Value = 3672.2727
$ws.Range("M79").Value = -1490.4333
$ws.Range("N79").Value = -5856.2727

$ws.Range("H86").Value = 1866.2142
$ws.Range("I86").Value = 1766.3334
$ws.Range("J86").Value = 2046
$ws.Range("K86").Value = 1766.3334
$ws.Range("L86").Value = 2046
$ws.Range("M86").Value = -643.3334
$ws.Range("N86").Value = -4292

$ws.Range("H89").Value = 1866.2142
$ws.Range("I89").Value = 1766.3334
$ws.Range("J89").Value = 2046
$ws.Range("K89").Value = 8831.666999999999
$ws.Range("L89").Value = 10230
$ws.Range("M89").Value = -3215.666999999999
$ws.Range("N89").Value = -21462

$ws.Range("H92").Value = 272.85184
$ws.Range("I92").Value = 244.88461
$ws.Range("K92").Value = 244.88461
$ws.Range("M92").Value = 1003.11539

$ws.Range("H98").Value = 1315.9474
$ws.Range("I98").Value = 1187.6875
$ws.Range("K98").Value = 1187.6875
$ws.Range("M98").Value = 310.3125

$ws.Range("H106").Value = 3642.875
$ws.Range("I106").Value = 2171.5
$ws.Range("J106").Value = 6095.1665
$ws.Range("K106").Value = 2171.5
$ws.Range("L106").Value = 6095.1665
$ws.Range("M106").Value = -1540.5
$ws.Range("N106").Value = -7357.1665

$ws.Range("H116").Value = 8552291
$ws.Range("I116").Value = 19232906
$ws.Range("J116").Value = 7800
$ws.Range("K116").Value = 19232906
$ws.Range("L116").Value = 7800
$ws.Range("M116").Value = -19229464
$ws.Range("N116").Value = -14684

$ws.Range("H122").Value = 1315.9474
$ws.Range("I122").Value = 1187.6875
$ws.Range("K122").Value = 3563.0625
$ws.Range("M122").Value = -1113.0625

$ws.Range("H125").Value = 5150
$ws.Range("J125").Value = 5150
$ws.Range("L125").Value = 46350
$ws.Range("N125").Value = -51270

$ws.Range("H132").Value = 2903.6
$ws.Range("I132").Value = 2389.347
$ws.Range("J132").Value = 7103.3335
$ws.Range("K132").Value = 7168.041000000001
$ws.Range("L132").Value = 21310.0005
$ws.Range("M132").Value = -4638.041000000001
$ws.Range("N132").Value = -26370.0005

$ws.Range("H137").Value = 21593.863
$ws.Range("I137").Value = 1205.1389
$ws.Range("J137").Value = 70526.8
$ws.Range("K137").Value = 3615.4167
$ws.Range("L137").Value = 211580.4
$ws.Range("M137").Value = -1065.4167
$ws.Range("N137").Value = -216680.4

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 2711311
$ws.Range("I74").Value = 1301.9412
$ws.Range("K74").Value = 1301.9412
$ws.Range("M74").Value = -427.9412

$ws.Range("H77").Value = 2711311
$ws.Range("I77").Value = 1301.9412
$ws.Range("K77").Value = 6509.706
$ws.Range("M77").Value = -2141.706

$ws.Range("H110").Value = 2892.818
$ws.Range("I110").Value = 2261.5715
$ws.Range("J110").Value = 3997.5
$ws.Range("K110").Value = 2261.5715
$ws.Range("L110").Value = 3997.5
$ws.Range("M110").Value = -216.5715
$ws.Range("N110").Value = -8087.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1847.9714
$ws.Range("I20").Value = 1256.909
$ws.Range("J20").Value = 2848.2307
$ws.Range("K20").Value = 1256.909
$ws.Range("L20").Value = 2848.2307
$ws.Range("M20").Value = -1009.909
$ws.Range("N20").Value = -3342.2307

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 8940.799999999999
$ws.Range("J2").Value = 19000
$ws.Range("L2").Value = 19000
$ws.Range("N2").Value = -19226

$ws.Range("H31").Value = 1881.92
$ws.Range("I31").Value = 1387.8462
$ws.Range("K31").Value = 1387.8462
$ws.Range("M31").Value = -1092.8462

$ws.Range("H34").Value = 1881.92
$ws.Range("I34").Value = 1387.8462
$ws.Range("K34").Value = 1387.8462
$ws.Range("M34").Value = -1185.8462

$ws.Range("H58").Value = 882.4074000000001
$ws.Range("I58").Value = 785.55554
$ws.Range("J58").Value = 1366.6666
$ws.Range("K58").Value = 785.55554
$ws.Range("L58").Value = 1366.6666
$ws.Range("M58").Value = -582.55554
$ws.Range("N58").Value = -1772.6666

$ws.Range("H132").Value = 1348.3829
$ws.Range("I132").Value = 1175.2683
$ws.Range("J132").Value = 2531.3333
$ws.Range("K132").Value = 3525.8049
$ws.Range("L132").Value = 7593.999899999999
$ws.Range("M132").Value = -995.8049000000001
$ws.Range("N132").Value = -12653.9999

$ws.Range("H136").Value = 882.4074000000001
$ws.Range("I136").Value = 785.55554
$ws.Range("J136").Value = 1366.6666
$ws.Range("K136").Value = 2356.66662
$ws.Range("L136").Value = 4099.9998
$ws.Range("M136").Value = 193.33338
$ws.Range("N136").Value = -9199.9998

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H63").Value = 2212.3333
$ws.Range("J63").Value = 3200
$ws.Range("L63").Value = 9600
$ws.Range("N63").Value = -11098

$ws.Range("H66").Value = 2212.3333
$ws.Range("J66").Value = 3200
$ws.Range("L66").Value = 28800
$ws.Range("N66").Value = -36288

$ws.Range("H107").Value = 296
$ws.Range("I107").Value = 196.66667
$ws.Range("J107").Value = 445
$ws.Range("K107").Value = 590.00001
$ws.Range("L107").Value = 1335
$ws.Range("M107").Value = 1329.99999
$ws.Range("N107").Value = -5175

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6422.2964
$ws.Range("I70").Value = 7136.4736
$ws.Range("J70").Value = 4726.125
$ws.Range("K70").Value = 7136.4736
$ws.Range("L70").Value = 4726.125
$ws.Range("M70").Value = -6866.4736
$ws.Range("N70").Value = -5266.125

$ws.Range("H73").Value = 6422.2964
$ws.Range("I73").Value = 7136.4736
$ws.Range("J73").Value = 4726.125
$ws.Range("K73").Value = 7136.4736
$ws.Range("L73").Value = 4726.125
$ws.Range("M73").Value = -6200.4736
$ws.Range("N73").Value = -6598.125

$ws.Range("H113").Value = 13447.556
$ws.Range("I113").Value = 1838
$ws.Range("J113").Value = 36666.668
$ws.Range("K113").Value = 1838
$ws.Range("L113").Value = 36666.668
$ws.Range("M113").Value = 332
$ws.Range("N113").Value = -41006.668

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2224
$ws.Range("I7").Value = 952
$ws.Range("J7").Value = 2860
$ws.Range("K7").Value = 952
$ws.Range("L7").Value = 2860
$ws.Range("M7").Value = -840
$ws.Range("N7").Value = -3084

$ws.Range("H61").Value = 2358.5
$ws.Range("I61").Value = 1316
$ws.Range("J61").Value = 2984
$ws.Range("K61").Value = 1316
$ws.Range("L61").Value = 2984
$ws.Range("M61").Value = -1114
$ws.Range("N61").Value = -3388

$ws.Range("H108").Value = 30000
$ws.Range("J108").Value = 30000
$ws.Range("L108").Value = 30000
$ws.Range("N108").Value = -37680

$ws.Range("H113").Value = 2358.5
$ws.Range("I113").Value = 1316
$ws.Range("J113").Value = 2984
$ws.Range("K113").Value = 1316
$ws.Range("L113").Value = 2984
$ws.Range("M113").Value = 854
$ws.Range("N113").Value = -7324

$ws.Range("H122").Value = 3378.9285
$ws.Range("I122").Value = 2868
$ws.Range("J122").Value = 3762.125
$ws.Range("K122").Value = 8604
$ws.Range("L122").Value = 11286.375
$ws.Range("M122").Value = -6154
$ws.Range("N122").Value = -16186.375

$ws.Range("H126").Value = 2224
$ws.Range("I126").Value = 952
$ws.Range("J126").Value = 2860
$ws.Range("K126").Value = 2856
$ws.Range("L126").Value = 8580
$ws.Range("M126").Value = -386
$ws.Range("N126").Value = -13520

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 934.7761
$ws.Range("I132").Value = 598.6964
$ws.Range("J132").Value = 2645.7273
$ws.Range("K132").Value = 1796.0892
$ws.Range("L132").Value = 7937.1819
$ws.Range("M132").Value = 733.9107999999999
$ws.Range("N132").Value = -12997.1819
